# Applies the "Handles float input without breaking stuff" marksheet update:
#  - Summary rows 10-12 get real computed scores (instead of placeholder
#    zeros / "Absent"), and their row-label cells (A10:A12) pick up the
#    "mtitleStyle" look already used by the row-9 header.
#  - C11's "-1" marking penalty was stored as text; it becomes a real
#    number so downstream SUM()s do not choke on it.
#  - The per-question "Student Ans" block only needs two answer columns
#    (A/B); the unused duplicate D/E and G/H blocks are cleared out from
#    row 19 down, and the G/H block is dropped everywhere (rows 15-18
#    too). Column A is populated with the student's chosen option,
#    colour-coded green/red/black for correct / incorrect / unattempted,
#    reusing the correctStyle/incorrectStyle/normalStyle cell styles that
#    already exist in the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# 1. Summary block (rows 10-12): real values instead of placeholders.
# ---------------------------------------------------------------------

# A10/A11/A12 ("No.", "Marking", "Total") adopt the same style as the
# other row labels in that block (e.g. A9, already "mtitleStyle").
$ws.Cells.Item(9, 1).Copy() | Out-Null
$ws.Range("A10:A12").PasteSpecial($xlPasteFormats) | Out-Null

$ws.Cells.Item(10, 2).Value = 19
$ws.Cells.Item(10, 3).Value = 3
$ws.Cells.Item(10, 4).Value = 6
$ws.Cells.Item(10, 5).Value = 28

$ws.Cells.Item(11, 2).Value = 4
# C11 used to hold the text "-1"; store it as a real number instead.
$ws.Cells.Item(11, 3).Value = -1

$ws.Cells.Item(12, 2).Value = 76
$ws.Cells.Item(12, 3).Value = -3
$ws.Cells.Item(12, 5).Value = "73/112"

# ---------------------------------------------------------------------
# 2. Drop the unused third "Student Ans / Correct Ans" block (G:H) for
#    every question row, and the duplicate D:E block from row 19 down.
# ---------------------------------------------------------------------

$ws.Range("G15:H40").Clear() | Out-Null
$ws.Range("D19:E40").Clear() | Out-Null

# ---------------------------------------------------------------------
# 3. Column A: fill in the student's answer per question, styled to
#    show correct (green), incorrect (red) or unattempted (black).
# ---------------------------------------------------------------------

# correctStyle (green) template already used by B10.
$ws.Cells.Item(10, 2).Copy() | Out-Null
$ws.Range("A16:A18").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A21:A27").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A29:A33").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A35:A35").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A38:A39").PasteSpecial($xlPasteFormats) | Out-Null
# D16 also becomes a correct-answer cell.
$ws.Range("D16").PasteSpecial($xlPasteFormats) | Out-Null

# incorrectStyle (red) template already used by C10.
$ws.Cells.Item(10, 3).Copy() | Out-Null
$ws.Range("A20:A20").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A37:A37").PasteSpecial($xlPasteFormats) | Out-Null
$ws.Range("A40:A40").PasteSpecial($xlPasteFormats) | Out-Null

# Rows 19, 28, 34, 36 remain unattempted (already normalStyle / blank) -
# nothing to change there.

$ws.Cells.Item(16, 1).Value = "Option A"
$ws.Cells.Item(17, 1).Value = "Option D"
$ws.Cells.Item(18, 1).Value = "Option B"
$ws.Cells.Item(20, 1).Value = "Option C"
$ws.Cells.Item(21, 1).Value = "Option C"
$ws.Cells.Item(22, 1).Value = "Option D"
$ws.Cells.Item(23, 1).Value = "Option D"
$ws.Cells.Item(24, 1).Value = "Option A"
$ws.Cells.Item(25, 1).Value = "Option A"
$ws.Cells.Item(26, 1).Value = "Option C"
$ws.Cells.Item(27, 1).Value = "Option A"
$ws.Cells.Item(29, 1).Value = "Option D"
$ws.Cells.Item(30, 1).Value = "Option B"
$ws.Cells.Item(31, 1).Value = "Option D"
$ws.Cells.Item(32, 1).Value = "Option C"
$ws.Cells.Item(33, 1).Value = "Option D"
$ws.Cells.Item(35, 1).Value = "Option D"
$ws.Cells.Item(37, 1).Value = "Option C"
$ws.Cells.Item(38, 1).Value = "Option A"
$ws.Cells.Item(39, 1).Value = "Option D"
$ws.Cells.Item(40, 1).Value = "Option B"

$ws.Cells.Item(16, 4).Value = "Option A"

Write-Host "Marksheet updated"
